$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 625, shifting rows 625:742 down to 627:744.
$ws.Rows.Item(625).Resize(2).Insert()

# New row 625 - "Primera" quality entry for the new date (44637 = 2022-03-17).
$ws.Cells.Item(625, 1).Value = 6
$ws.Cells.Item(625, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(625, 3).Value = "Metropolitana"
$ws.Cells.Item(625, 4).Value = 44637
$ws.Cells.Item(625, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(625, 5).Value = 13
$ws.Cells.Item(625, 6).Value = 100112009
$ws.Cells.Item(625, 7).Value = "Acelga"
$ws.Cells.Item(625, 8).Value = "Sin especificar"
$ws.Cells.Item(625, 9).Value = "Primera"
$ws.Cells.Item(625, 10).Value = 230
$ws.Cells.Item(625, 11).Value = 12000
$ws.Cells.Item(625, 12).Value = 12000
$ws.Cells.Item(625, 13).Value = 12000
$ws.Cells.Item(625, 14).Value = "`$/docena de atados"
$ws.Cells.Item(625, 15).Value = "Región Metropolitana"
$ws.Cells.Item(625, 16).Value = 4000
$ws.Cells.Item(625, 17).Value = 3
$ws.Cells.Item(625, 18).Value = "Hortaliza"

# New row 626 - "Segunda" quality entry for the new date (44637 = 2022-03-17).
$ws.Cells.Item(626, 1).Value = 6
$ws.Cells.Item(626, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(626, 3).Value = "Metropolitana"
$ws.Cells.Item(626, 4).Value = 44637
$ws.Cells.Item(626, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(626, 5).Value = 13
$ws.Cells.Item(626, 6).Value = 100112009
$ws.Cells.Item(626, 7).Value = "Acelga"
$ws.Cells.Item(626, 8).Value = "Sin especificar"
$ws.Cells.Item(626, 9).Value = "Segunda"
$ws.Cells.Item(626, 10).Value = 150
$ws.Cells.Item(626, 11).Value = 10000
$ws.Cells.Item(626, 12).Value = 10000
$ws.Cells.Item(626, 13).Value = 10000
$ws.Cells.Item(626, 14).Value = "`$/docena de atados"
$ws.Cells.Item(626, 15).Value = "Región Metropolitana"
$ws.Cells.Item(626, 16).Value = 3333
$ws.Cells.Item(626, 17).Value = 3
$ws.Cells.Item(626, 18).Value = "Hortaliza"
